# The wml.xsd schema expects run-property child elements in a fixed
# sequence (rFonts, b, bc, i, ... , color, ...). Several of the Pandoc
# "Tok" character styles in styles.xml had <w:color/> emitted before
# <w:b/>/<w:i/>, which OOXMLValidatorCLI flags as a schema violation
# even though xmllint stays quiet. Re-assert the same Bold/Italic/Color
# values on each affected style so the engine re-emits <w:rPr> in the
# schema-correct order (bold/italic before color) without changing any
# actual formatting value.

$d = $word.ActiveDocument
$styles = $d.Styles

$boldColorStyles = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
foreach ($name in $boldColorStyles) {
    $s = $styles.Item($name)
    $s.Font.Bold = $true
}

$italicColorStyles = @("CommentTok", "DocumentationTok")
foreach ($name in $italicColorStyles) {
    $s = $styles.Item($name)
    $s.Font.Italic = $true
}

$boldItalicColorStyles = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")
foreach ($name in $boldItalicColorStyles) {
    $s = $styles.Item($name)
    $s.Font.Bold = $true
    $s.Font.Italic = $true
}
